$wb = $excel.ActiveWorkbook

$wsAfstand = $wb.Worksheets.Item("Afstand_km")
$wsDuur = $wb.Worksheets.Item("Duur_minuten")

# --- Afstand_km (sheet1) ---
$wsAfstand.Range("D6").Value = 135.64
$wsAfstand.Range("N6").Value = 126.38
$wsAfstand.Range("D7").Value = 103.74
$wsAfstand.Range("N7").Value = 94.47
$wsAfstand.Range("D10").Value = 73.81
$wsAfstand.Range("I10").Value = 254.74
$wsAfstand.Range("N10").Value = 64.54000000000001
$wsAfstand.Range("D11").Value = 140.88
$wsAfstand.Range("N11").Value = 131.62
$wsAfstand.Range("C13").Value = 218.8
$wsAfstand.Range("D13").Value = 242.26
$wsAfstand.Range("I13").Value = 399.42
$wsAfstand.Range("N13").Value = 232.99

# --- Duur_minuten (sheet2) ---
$wsDuur.Range("I2").Value = 8165
$wsDuur.Range("B3").Value = 4105.97
$wsDuur.Range("D3").Value = 1640.87
$wsDuur.Range("I3").Value = 10752.52
$wsDuur.Range("L3").Value = 9150.639999999999
$wsDuur.Range("N3").Value = 1119.72
$wsDuur.Range("B4").Value = 3281.09
$wsDuur.Range("I4").Value = 10094.5
$wsDuur.Range("B5").Value = 2543.26
$wsDuur.Range("I5").Value = 8668.26
$wsDuur.Range("B6").Value = 6054.75
$wsDuur.Range("D6").Value = 6268.44
$wsDuur.Range("N6").Value = 5747.29
$wsDuur.Range("B7").Value = 5740.73
$wsDuur.Range("C7").Value = 4072.27
$wsDuur.Range("D7").Value = 4802.14
$wsDuur.Range("F7").Value = 2769.91
$wsDuur.Range("I7").Value = 11591.67
$wsDuur.Range("J7").Value = 1851.65
$wsDuur.Range("M7").Value = 5818.54
$wsDuur.Range("N7").Value = 4280.99
$wsDuur.Range("B8").Value = 2067.51
$wsDuur.Range("E8").Value = 2459.81
$wsDuur.Range("I8").Value = 7036.3
$wsDuur.Range("B9").Value = 8030.14
$wsDuur.Range("E9").Value = 8461.629999999999
$wsDuur.Range("F9").Value = 11311.15
$wsDuur.Range("G9").Value = 11311.46
$wsDuur.Range("K9").Value = 9171.120000000001
$wsDuur.Range("M9").Value = 13321.27
$wsDuur.Range("B10").Value = 5397
$wsDuur.Range("D10").Value = 3739.54
$wsDuur.Range("I10").Value = 11247.93
$wsDuur.Range("K10").Value = 4874.26
$wsDuur.Range("N10").Value = 3218.39
$wsDuur.Range("B11").Value = 5930.4
$wsDuur.Range("D11").Value = 6465.19
$wsDuur.Range("G11").Value = 4028.84
$wsDuur.Range("J11").Value = 4931.44
$wsDuur.Range("N11").Value = 5944.03
$wsDuur.Range("B12").Value = 6433.21
$wsDuur.Range("E12").Value = 7461.29
$wsDuur.Range("F12").Value = 11001.36
$wsDuur.Range("G12").Value = 10455.67
$wsDuur.Range("H12").Value = 5876.22
$wsDuur.Range("I12").Value = 3743.44
$wsDuur.Range("J12").Value = 9976.68
$wsDuur.Range("K12").Value = 8861.33
$wsDuur.Range("M12").Value = 13944.51
$wsDuur.Range("B13").Value = 9406.82
$wsDuur.Range("C13").Value = 8740.139999999999
$wsDuur.Range("D13").Value = 9470.01
$wsDuur.Range("I13").Value = 13156.62
$wsDuur.Range("N13").Value = 8948.85
$wsDuur.Range("B14").Value = 3532.92
$wsDuur.Range("I14").Value = 10346.33
